$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 with new values
$ws.Cells.Item(2, 2).Value = 606
$ws.Cells.Item(2, 3).Value = 0.03254204303681323
$ws.Cells.Item(2, 5).Value = 0.9999016790546436
$ws.Cells.Item(2, 6).Value = 0.09977054595947266
$ws.Cells.Item(2, 7).Value = 0.7742536362831121

$ws.Cells.Item(3, 2).Value = 5673
$ws.Cells.Item(3, 3).Value = 0.01789872786270224
$ws.Cells.Item(3, 5).Value = 0.09998305469224604
$ws.Cells.Item(3, 6).Value = 0.7767543792724609
$ws.Cells.Item(3, 7).Value = 0.4428740998807554

# Row 4
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(4, 1).PasteSpecial(-4122)
$ws.Cells.Item(4, 2).Value = 13767
$ws.Cells.Item(4, 3).Value = 0.01643449281948481
$ws.Cells.Item(4, 4).Value = 0.01
$ws.Cells.Item(4, 5).Value = 0.009997121167773834
$ws.Cells.Item(4, 6).Value = 1.899664402008057
$ws.Cells.Item(4, 7).Value = 0.3368040026290677

# Row 5
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)
$ws.Cells.Item(5, 2).Value = 20622
$ws.Cells.Item(5, 3).Value = 0.01628805934426549
$ws.Cells.Item(5, 4).Value = 0.001
$ws.Cells.Item(5, 5).Value = 0.0009979150444940192
$ws.Cells.Item(5, 6).Value = 2.787895202636719
$ws.Cells.Item(5, 7).Value = 0.3163171269278453

# Row 6
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(6, 1).PasteSpecial(-4122)
$ws.Cells.Item(6, 2).Value = 21451
$ws.Cells.Item(6, 3).Value = 0.01627343664101082
$ws.Cells.Item(6, 4).Value = 0.0001
$ws.Cells.Item(6, 5).Value = 0.00009926314491103337
$ws.Cells.Item(6, 6).Value = 2.893169164657593
$ws.Cells.Item(6, 7).Value = 0.3147820630490531

# Row 7
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(7, 1).PasteSpecial(-4122)
$ws.Cells.Item(7, 2).Value = 21504
$ws.Cells.Item(7, 3).Value = 0.01627195917621057
$ws.Cells.Item(7, 4).Value = 0.00001
$ws.Cells.Item(7, 5).Value = 0.000008464164299558957
$ws.Cells.Item(7, 6).Value = 2.900102138519287
$ws.Cells.Item(7, 7).Value = 0.31469217492118

# Row 8
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122)
$ws.Cells.Item(8, 2).Value = 22051
$ws.Cells.Item(8, 3).Value = 0.01627181648572707
$ws.Cells.Item(8, 4).Value = 0.000001
$ws.Cells.Item(8, 5).Value = 0.0000003050127683218948
$ws.Cells.Item(8, 6).Value = 2.969652414321899
$ws.Cells.Item(8, 7).Value = 0.3143282789609663

# Row 9
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4122)
$ws.Cells.Item(9, 2).Value = 24802
$ws.Cells.Item(9, 3).Value = 0.01627182112465549
$ws.Cells.Item(9, 4).Value = 0.0000001
$ws.Cells.Item(9, 5).Value = 0.00000001992308490776627
$ws.Cells.Item(9, 6).Value = 3.332269430160522
$ws.Cells.Item(9, 7).Value = 0.3143209987791079

$excel.CutCopyMode = $false
